$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "Leer parte de las tareas"
$ws.Range("C8").Value = "Nada"
$ws.Range("C9").Value = "Nada"

$ws.Range("D7").Value = "Nada"
$ws.Range("D8").Value = "Resumir las lecturas sobre el SEMAT, mirar que es IGLUW"
$ws.Range("D9").Value = "Nada"

$ws.Range("E7").Value = "Se logro lo propuesto"
$ws.Range("E8").Value = "Asistir y aportar a la reunion para practicar la exposicion de lo que necesitamos para la clase de mañana"
$ws.Range("E9").Value = "Nada"

$ws.Range("F13").Select() | Out-Null
